# Remove column M from the alcohol measurement data on Sheet1.
# The column to its right (old N) shifts left to become the new M,
# and the used range shrinks from A1:N119 to A1:M119.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

[void]$ws.Range("M1:M119").Delete()
[void]$ws.Range("M1").Select()
